$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, pushing existing rows 172:185 down to 173:186
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with the new price-report record
$ws.Range("A172").Value = 4
$ws.Range("B172").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C172").Value = "Los Lagos"
$ws.Range("D172").Value = 45265
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 100112022
$ws.Range("G172").Value = "Arveja Verde"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 100
$ws.Range("K172").Value = 27000
$ws.Range("L172").Value = 27000
$ws.Range("M172").Value = 27000
$ws.Range("N172").Value = "$/saco 25 kilos"
$ws.Range("O172").Value = "Región del Maule"
$ws.Range("P172").Value = 1080
$ws.Range("Q172").Value = 25
$ws.Range("R172").Value = "Hortaliza"
